# Automatische test-sync: 2025-07-22 12:32:50
# Appends Testmail #10 ("Kun je dit inkopen voor ons project?") to the
# "Logs" sheet as row 10, extends the conditional-formatting ranges that
# covered rows 2:9 to 2:10, and bumps the "Productinformatie" tally on the
# "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new log row (row 10) -----------------------------------

$antwoord = @"
Beste afzender,
Dank voor uw e-mail. Om uw verzoek te kunnen verwerken, zou ik graag meer informatie ontvangen over het product dat u wilt inkopen voor uw project. Kunt u meer details geven over het product en de hoeveelheid die u nodig heeft? Eventueel een offerte of specificaties kunnen helpen ons te helpen om dit voor u te regelen.
Met vriendelijke groet,
[Naam]   
E-mailassistent
"@

$logs.Range("A10").Value = "Kun je dit inkopen voor ons project?"
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Testmail #10: Kun je dit inkopen voor ons project?"
$logs.Range("D10").Value = "Productinformatie"
$logs.Range("E10").Value = $antwoord
$logs.Range("F10").Value = "2025-07-22 12:32:26"
$logs.Range("G10").Value = "Ja"
$logs.Range("H10").Value = "Nee"
$logs.Range("I10").Value = "Ja"
$logs.Range("J10").Value = "Ja"

# --- 2. Extend the conditional formatting ranges D2:D9 .. J2:J9 to ..:10 --

$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range("$($col)2:$($col)9")
    $newRange = $logs.Range("$($col)2:$($col)10")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the Dashboard "Productinformatie" count (3 -> 4) -----------

$dashboard.Range("B2").Value = 4
